# Update weekly Pepino dulce price-report rows (Hortaliza, Agricola del Norte S.A. de Arica).
# Values below reconstruct the new week of rows by shifting the reported
# observations forward one slot (row 2's old data is replaced, row 9/16/17 keep theirs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44377
$ws.Range('H2').Value = 'Cultivar IV Región'
$ws.Range('I2').Value = 'Primera'
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17600
$ws.Range('N2').Value = '$/bandeja 18 kilos'
$ws.Range('O2').Value = 'Provincia de Limarí'
$ws.Range("P2").Value = 978
$ws.Range("Q2").Value = 18

# Row 3
$ws.Range("D3").Value = 44435

# Row 4
$ws.Range("D4").Value = 44435

# Row 5
$ws.Range("D5").Value = 44363
$ws.Range('H5').Value = 'Cultivar IV Región'
$ws.Range("J5").Value = 140
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range('N5').Value = '$/bandeja 18 kilos'
$ws.Range('O5').Value = 'Provincia de Limarí'
$ws.Range("P5").Value = 806
$ws.Range("Q5").Value = 18

# Row 6
$ws.Range("D6").Value = 44391
$ws.Range('H6').Value = 'Cultivar IV Región'
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range('N6').Value = '$/bandeja 18 kilos'
$ws.Range('O6').Value = 'Provincia de Limarí'
$ws.Range("P6").Value = 861
$ws.Range("Q6").Value = 18

# Row 7
$ws.Range("D7").Value = 44398
$ws.Range('H7').Value = 'Cultivar IV Región'
$ws.Range('I7').Value = 'Primera'
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17500
$ws.Range('N7').Value = '$/bandeja 18 kilos'
$ws.Range('O7').Value = 'Provincia de Limarí'
$ws.Range("P7").Value = 972
$ws.Range("Q7").Value = 18

# Row 8
$ws.Range("D8").Value = 44398
$ws.Range('I8').Value = 'Segunda'
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("P8").Value = 861

# Row 10
$ws.Range("D10").Value = 44526
$ws.Range('H10').Value = 'Cultivar XV región'
$ws.Range('I10').Value = 'Primera'
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 5500
$ws.Range("M10").Value = 5250
$ws.Range('N10').Value = '$/caja 10 kilos'
$ws.Range('O10').Value = 'Región de Arica y Parinacota'
$ws.Range("P10").Value = 525
$ws.Range("Q10").Value = 10

# Row 11
$ws.Range("D11").Value = 44526
$ws.Range('H11').Value = 'Cultivar XV región'
$ws.Range('I11').Value = 'Segunda'
$ws.Range("K11").Value = 4000
$ws.Range("L11").Value = 4500
$ws.Range("M11").Value = 4250
$ws.Range('N11').Value = '$/caja 10 kilos'
$ws.Range('O11').Value = 'Región de Arica y Parinacota'
$ws.Range("P11").Value = 425
$ws.Range("Q11").Value = 10

# Row 12
$ws.Range("D12").Value = 44526
$ws.Range('I12').Value = 'Tercera'
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3500
$ws.Range("M12").Value = 3250
$ws.Range("P12").Value = 325

# Row 13
$ws.Range("D13").Value = 44554
$ws.Range('H13').Value = 'Cultivar XV región'
$ws.Range('I13').Value = 'Primera'
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = 5500
$ws.Range('N13').Value = '$/caja 10 kilos'
$ws.Range('O13').Value = 'Región de Arica y Parinacota'
$ws.Range("P13").Value = 550
$ws.Range("Q13").Value = 10

# Row 14
$ws.Range("D14").Value = 44405
$ws.Range("J14").Value = 140

# Row 15
$ws.Range("D15").Value = 44221
$ws.Range('H15').Value = 'Cultivar XV región'
$ws.Range('I15').Value = 'Primera'
$ws.Range("J15").Value = 140
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 6000
$ws.Range("M15").Value = 5500
$ws.Range('N15').Value = '$/caja 10 kilos'
$ws.Range('O15').Value = 'Región de Arica y Parinacota'
$ws.Range("P15").Value = 550
$ws.Range("Q15").Value = 10

# Row 18
$ws.Range("D18").Value = 44433
$ws.Range('I18').Value = 'Segunda'
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17500
$ws.Range("P18").Value = 972

# Row 19
$ws.Range("D19").Value = 44433
$ws.Range('I19').Value = 'Tercera'
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 14500
$ws.Range("P19").Value = 806

# Row 20
$ws.Range("D20").Value = 44211
$ws.Range('H20').Value = 'Cultivar XV región'
$ws.Range("J20").Value = 140
$ws.Range("K20").Value = 4500
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = 4750
$ws.Range('N20').Value = '$/caja 10 kilos'
$ws.Range('O20').Value = 'Región de Arica y Parinacota'
$ws.Range("P20").Value = 475
$ws.Range("Q20").Value = 10

# Row 21
$ws.Range("D21").Value = 44454
$ws.Range('H21').Value = 'Cultivar IV Región'
$ws.Range("J21").Value = 160
$ws.Range("K21").Value = 19000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 19500
$ws.Range('N21').Value = '$/bandeja 18 kilos'
$ws.Range('O21').Value = 'Provincia de Limarí'
$ws.Range("P21").Value = 1083
$ws.Range("Q21").Value = 18
